{"js": "// Remove the \"Appendix: Quick prototype\" section: its Heading2 paragraph,\n// the figure captions, and the embedded screenshots that follow it, up to\n// (but not including) the next \"Appendix: Links\" Heading2 paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the start of the section to remove: the Heading2 paragraph whose\n// text is \"Appendix: Quick prototype\".\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\" && items[i].text.trim() === \"Appendix: Quick prototype\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Locate the end of the section: the next Heading2 paragraph after\n  // startIndex (exclusive) \u2014 everything in between gets removed.\n  let endIndex = items.length;\n  for (let i = startIndex + 1; i < items.length; i++) {\n    if (items[i].style === \"Heading 2\") {\n      endIndex = i;\n      break;\n    }\n  }\n\n  // Delete in reverse order so earlier indices stay valid as we go.\n  for (let i = endIndex - 1; i >= startIndex; i--) {\n    items[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"Appendix: Quick prototype\" section: its Heading2 paragraph,\n# the figure captions, and the embedded screenshots that follow it, up to\n# (but not including) the next \"Appendix: Links\" Heading2 paragraph.\n$d = $word.ActiveDocument\n\n# Locate the start of the section: the paragraph containing\n# \"Appendix: Quick prototype\".\n$startRange = $d.Content\n[void]$startRange.Find.Execute(\"Appendix: Quick prototype\")\n\nif ($startRange.Find.Found) {\n    $startPara = $startRange.Paragraphs.First\n    $startPos = $startPara.Range.Start\n\n    # Locate the end of the section: the next \"Appendix: Links\" heading\n    # paragraph that follows the start point.\n    $searchRange = $d.Range($startRange.End, $d.Content.End)\n    [void]$searchRange.Find.Execute(\"Appendix: Links\")\n\n    if ($searchRange.Find.Found) {\n        $endPara = $searchRange.Paragraphs.First\n        $endPos = $endPara.Range.Start\n    } else {\n        $endPos = $d.Content.End\n    }\n\n    $deleteRange = $d.Range($startPos, $endPos)\n    $deleteRange.Delete()\n}\n"}
